# Applies the data corrections described in the commit:
#   "changed MP time limit and corrected error in fixed recourse data"
#
# Sheet1 holds the per-instance summary (objective, solve time, num cuts,
# num variables, num cons, num quad cons). Sheets "1".."10" hold the
# per-iteration master-problem (MP) detail for each instance.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1 : summary table (rows 2-11, columns B,C,F,G,H,I)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$summary = @(
    @{ Row = 2;  B = -105.20922805100011; C = 9.36610749 },
    @{ Row = 3;  B = -100.25578091436573; C = 1.417407134 },
    @{ Row = 4;  B = -103.67140171307204; C = 1.508649313 },
    @{ Row = 5;  B = -103.09145004118963; C = 1.509030697 },
    @{ Row = 6;  B = -102.05315976383038; C = 1.338546347 },
    @{ Row = 7;  B = -102.01057733765836; C = 1.106352739 },
    @{ Row = 8;  B = -97.44343294495316;  C = 0.993129861 },
    @{ Row = 9;  B = -102.53029754612697; C = 1.720042824 },
    @{ Row = 10; B = -102.29121933446089; C = 1.269073295 },
    @{ Row = 11; B = -99.01915297564395;  C = 0.931378026 }
)

foreach ($entry in $summary) {
    $r = $entry.Row
    $ws1.Cells.Item($r, 2).Value = $entry.B   # B: objective
    $ws1.Cells.Item($r, 3).Value = $entry.C   # C: solve time
    $ws1.Cells.Item($r, 6).Value = 20         # F: num cuts
    $ws1.Cells.Item($r, 7).Value = 4490       # G: num variables
    $ws1.Cells.Item($r, 8).Value = 4900       # H: num cons
    $ws1.Cells.Item($r, 9).Value = 400        # I: num quad cons
}

# ---------------------------------------------------------------------
# Sheets "1".."10" : per-iteration MP detail (rows 2-3)
# ---------------------------------------------------------------------

# Sheet "1"
$ws = $wb.Worksheets.Item("1")
$ws.Cells.Item(2, 4).Value = 0.8300394577015381
$ws.Cells.Item(2, 5).Value = 49.3648
$ws.Cells.Item(3, 2).Value = -105.20922805100011
$ws.Cells.Item(3, 3).Value = 0.017387750492357718
$ws.Cells.Item(3, 4).Value = 1.1209461137015382

# Sheet "2"
$ws = $wb.Worksheets.Item("2")
$ws.Cells.Item(2, 4).Value = 0.011858894009643555
$ws.Cells.Item(2, 5).Value = 51.44005
$ws.Cells.Item(3, 2).Value = -100.25578091436573
$ws.Cells.Item(3, 3).Value = 0.0
$ws.Cells.Item(3, 4).Value = 1.2577127725305175

# Sheet "3"
$ws = $wb.Worksheets.Item("3")
$ws.Cells.Item(2, 4).Value = 0.022247979205322267
$ws.Cells.Item(2, 5).Value = 48.49419
$ws.Cells.Item(3, 2).Value = -103.67140171307204
$ws.Cells.Item(3, 4).Value = 1.3519543067456055

# Sheet "4"
$ws = $wb.Worksheets.Item("4")
$ws.Cells.Item(2, 4).Value = 0.027721889509155273
$ws.Cells.Item(2, 5).Value = 51.58398
$ws.Cells.Item(3, 2).Value = -103.09145004118963
$ws.Cells.Item(3, 4).Value = 1.3256724619320068

# Sheet "5"
$ws = $wb.Worksheets.Item("5")
$ws.Cells.Item(2, 4).Value = 0.025909265294921876
$ws.Cells.Item(2, 5).Value = 53.26468
$ws.Cells.Item(3, 2).Value = -102.05315976383038
$ws.Cells.Item(3, 3).Value = 0.0
$ws.Cells.Item(3, 4).Value = 1.1496742596733398

# Sheet "6"
$ws = $wb.Worksheets.Item("6")
$ws.Cells.Item(2, 4).Value = 0.03059385074194336
$ws.Cells.Item(2, 5).Value = 52.09259
$ws.Cells.Item(3, 2).Value = -102.01057733765836
$ws.Cells.Item(3, 3).Value = 0.0
$ws.Cells.Item(3, 4).Value = 0.9254119965684815

# Sheet "7"
$ws = $wb.Worksheets.Item("7")
$ws.Cells.Item(2, 4).Value = 0.03835181897802734
$ws.Cells.Item(2, 5).Value = 46.78724
$ws.Cells.Item(3, 2).Value = -97.44343294495316
$ws.Cells.Item(3, 3).Value = 0.0
$ws.Cells.Item(3, 4).Value = 0.8507319640587159

# Sheet "8"
$ws = $wb.Worksheets.Item("8")
$ws.Cells.Item(2, 4).Value = 0.04183498033972168
$ws.Cells.Item(2, 5).Value = 52.74845
$ws.Cells.Item(3, 2).Value = -102.53029754612697
$ws.Cells.Item(3, 4).Value = 1.4942494058953857

# Sheet "9"
$ws = $wb.Worksheets.Item("9")
$ws.Cells.Item(2, 4).Value = 0.013725472433715821
$ws.Cells.Item(2, 5).Value = 52.47722
$ws.Cells.Item(3, 2).Value = -102.29121933446089
$ws.Cells.Item(3, 3).Value = 0.09487514189111497
$ws.Cells.Item(3, 4).Value = 1.0940487608913574

# Sheet "10"
$ws = $wb.Worksheets.Item("10")
$ws.Cells.Item(2, 4).Value = 0.011581720487304688
$ws.Cells.Item(2, 5).Value = 49.97122
$ws.Cells.Item(3, 2).Value = -99.01915297564395
$ws.Cells.Item(3, 4).Value = 0.7937238566665039

Write-Host "Applied fixed-recourse data corrections to all sheets."
